$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 708.26666
$ws.Range("I33").Value = 560.4167
$ws.Range("J33").Value = 1299.6666
$ws.Range("K33").Value = 560.4167
$ws.Range("L33").Value = 1299.6666
$ws.Range("M33").Value = -331.4167
$ws.Range("N33").Value = -1757.6666
$ws.Range("H40").Value = 2817.0588
$ws.Range("I40").Value = 2222.5
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2222.5
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2047.5
$ws.Range("N40").Value = -3350
$ws.Range("H112").Value = 1654.3529
$ws.Range("I112").Value = 937.2857
$ws.Range("J112").Value = 2156.3
$ws.Range("K112").Value = 2811.8571
$ws.Range("L112").Value = 6468.900000000001
$ws.Range("M112").Value = -1703.8571
$ws.Range("N112").Value = -8684.900000000001
$ws.Range("H137").Value = 1617581.9
$ws.Range("I137").Value = 1158820.9
$ws.Range("J137").Value = 1923422.5
$ws.Range("K137").Value = 3476462.7
$ws.Range("L137").Value = 5770267.5
$ws.Range("M137").Value = -3473912.7
$ws.Range("N137").Value = -5775367.5
$ws.Range("H138").Value = 4253.4688
$ws.Range("J138").Value = 4046.889
$ws.Range("L138").Value = 12140.667
$ws.Range("N138").Value = -22420.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4133.467
$ws.Range("I2").Value = 4750.1113
$ws.Range("J2").Value = 3208.5
$ws.Range("K2").Value = 4750.1113
$ws.Range("L2").Value = 3208.5
$ws.Range("M2").Value = -4637.1113
$ws.Range("N2").Value = -3434.5
$ws.Range("H32").Value = 6181284
$ws.Range("I32").Value = 7100177.5
$ws.Range("J32").Value = 11569.286
$ws.Range("K32").Value = 7100177.5
$ws.Range("L32").Value = 11569.286
$ws.Range("M32").Value = -7099890.5
$ws.Range("N32").Value = -12143.286
$ws.Range("H61").Value = 581939.1
$ws.Range("I61").Value = 958477
$ws.Range("K61").Value = 958477
$ws.Range("M61").Value = -958265
$ws.Range("H102").Value = 4287.8335
$ws.Range("I102").Value = 3577.75
$ws.Range("K102").Value = 3577.75
$ws.Range("M102").Value = -1955.75
$ws.Range("H116").Value = 4133.467
$ws.Range("I116").Value = 4750.1113
$ws.Range("J116").Value = 3208.5
$ws.Range("K116").Value = 4750.1113
$ws.Range("L116").Value = 3208.5
$ws.Range("M116").Value = -2456.1113
$ws.Range("N116").Value = -7796.5
$ws.Range("H122").Value = 1522
$ws.Range("I122").Value = 1522
$ws.Range("K122").Value = 4566
$ws.Range("M122").Value = -2116
$ws.Range("H132").Value = 241782.9
$ws.Range("I132").Value = 404340.2
$ws.Range("K132").Value = 1213020.6
$ws.Range("M132").Value = -1210490.6
$ws.Range("H136").Value = 581939.1
$ws.Range("I136").Value = 958477
$ws.Range("K136").Value = 2875431
$ws.Range("M136").Value = -2872881

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4133.467
$ws.Range("I3").Value = 4750.1113
$ws.Range("J3").Value = 3208.5
$ws.Range("K3").Value = 4750.1113
$ws.Range("L3").Value = 3208.5
$ws.Range("M3").Value = -4636.1113
$ws.Range("N3").Value = -3436.5
$ws.Range("H20").Value = 2832.1333
$ws.Range("J20").Value = 3394.5833
$ws.Range("L20").Value = 3394.5833
$ws.Range("N20").Value = -3888.5833
$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877
$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384
$ws.Range("H94").Value = 738.9032
$ws.Range("I94").Value = 607.6667
$ws.Range("K94").Value = 607.6667
$ws.Range("M94").Value = -156.6667
$ws.Range("H107").Value = 2094.9092
$ws.Range("I107").Value = 2551.1333
$ws.Range("K107").Value = 2551.1333
$ws.Range("M107").Value = -631.1333
$ws.Range("H132").Value = 95389.5
$ws.Range("J132").Value = 95389.5
$ws.Range("L132").Value = 95389.5
$ws.Range("N132").Value = -105509.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7743.8335
$ws.Range("J31").Value = 10968
$ws.Range("L31").Value = 10968
$ws.Range("N31").Value = -11558
$ws.Range("H34").Value = 7743.8335
$ws.Range("J34").Value = 10968
$ws.Range("L34").Value = 10968
$ws.Range("N34").Value = -11372
$ws.Range("H58").Value = 459804.6
$ws.Range("I58").Value = 477104.78
$ws.Range("K58").Value = 477104.78
$ws.Range("M58").Value = -476901.78
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("N83").Value = 0
$ws.Range("H86").Value = 2453.8
$ws.Range("J86").Value = 2700.6
$ws.Range("L86").Value = 2700.6
$ws.Range("N86").Value = -4946.6
$ws.Range("H89").Value = 2453.8
$ws.Range("J89").Value = 2700.6
$ws.Range("L89").Value = 13503
$ws.Range("N89").Value = -24735
$ws.Range("H96").Value = 33494.5
$ws.Range("J96").Value = 33494.5
$ws.Range("L96").Value = 33494.5
$ws.Range("N96").Value = -38986.5
$ws.Range("H132").Value = 12633.424
$ws.Range("I132").Value = 13473.5
$ws.Range("K132").Value = 40420.5
$ws.Range("M132").Value = -37890.5
$ws.Range("H136").Value = 459804.6
$ws.Range("I136").Value = 477104.78
$ws.Range("K136").Value = 1431314.34
$ws.Range("M136").Value = -1428764.34
$ws.Range("N80").ClearContents()
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200032.64
$ws.Range("I2").Value = 333366.88
$ws.Range("K2").Value = 2000201.28
$ws.Range("M2").Value = -2000088.28
$ws.Range("H17").Value = 81.52381
$ws.Range("I17").Value = 36.473682
$ws.Range("J17").Value = 509.5
$ws.Range("K17").Value = 109.421046
$ws.Range("L17").Value = 1528.5
$ws.Range("M17").Value = 59.57895400000001
$ws.Range("N17").Value = -1866.5
$ws.Range("H34").Value = 1249.8334
$ws.Range("I34").Value = 874.75
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 2624.25
$ws.Range("L34").Value = 6000
$ws.Range("M34").Value = -2540.25
$ws.Range("N34").Value = -6168
$ws.Range("H55").Value = 2749.75
$ws.Range("I55").Value = 2999.5
$ws.Range("K55").Value = 8998.5
$ws.Range("M55").Value = -8821.5
$ws.Range("H113").Value = 1465.4736
$ws.Range("I113").Value = 1246.8
$ws.Range("J113").Value = 1498.6061
$ws.Range("K113").Value = 3740.4
$ws.Range("L113").Value = 4495.8183
$ws.Range("M113").Value = -1570.4
$ws.Range("N113").Value = -8835.818299999999
$ws.Range("H121").Value = 706.4167
$ws.Range("I121").Value = 698.7
$ws.Range("J121").Value = 711.9286
$ws.Range("K121").Value = 2096.1
$ws.Range("L121").Value = 2135.7858
$ws.Range("M121").Value = -786.1000000000004
$ws.Range("N121").Value = -4755.7858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4416.4165
$ws.Range("I70").Value = 4500.2
$ws.Range("K70").Value = 4500.2
$ws.Range("M70").Value = -4230.2
$ws.Range("H73").Value = 4416.4165
$ws.Range("I73").Value = 4500.2
$ws.Range("K73").Value = 4500.2
$ws.Range("M73").Value = -3564.2
$ws.Range("H80").Value = 164053.06
$ws.Range("I80").Value = 194930.66
$ws.Range("K80").Value = 194930.66
$ws.Range("M80").Value = -193932.66
$ws.Range("H83").Value = 164053.06
$ws.Range("I83").Value = 194930.66
$ws.Range("K83").Value = 974653.3
$ws.Range("M83").Value = -969661.3
$ws.Range("H97").Value = 1002.1177
$ws.Range("I97").Value = 740.7222
$ws.Range("J97").Value = 1296.1875
$ws.Range("K97").Value = 740.7222
$ws.Range("L97").Value = 1296.1875
$ws.Range("M97").Value = -244.7222
$ws.Range("N97").Value = -2288.1875
$ws.Range("H122").Value = 3886.1724
$ws.Range("I122").Value = 1946
$ws.Range("K122").Value = 5838
$ws.Range("M122").Value = -3388
$ws.Range("H126").Value = 697784.5
$ws.Range("I126").Value = 927713.3
$ws.Range("J126").Value = 7998
$ws.Range("K126").Value = 2783139.9
$ws.Range("L126").Value = 23994
$ws.Range("M126").Value = -2780669.9
$ws.Range("N126").Value = -28934
$ws.Range("H132").Value = 242764.4
$ws.Range("I132").Value = 345897.5
$ws.Range("K132").Value = 1037692.5
$ws.Range("M132").Value = -1035162.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3442.5
$ws.Range("I40").Value = 3442.5
$ws.Range("K40").Value = 3442.5
$ws.Range("M40").Value = -3306.5
$ws.Range("H132").Value = 655878.4399999999
$ws.Range("I132").Value = 912620.3
$ws.Range("J132").Value = 5465.533
$ws.Range("K132").Value = 2737860.9
$ws.Range("L132").Value = 16396.599
$ws.Range("M132").Value = -2735330.9
$ws.Range("N132").Value = -21456.599

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6944045
$ws.Range("I132").Value = 11842235
$ws.Range("J132").Value = 4941.5
$ws.Range("K132").Value = 35526705
$ws.Range("L132").Value = 14824.5
$ws.Range("M132").Value = -35524175
$ws.Range("N132").Value = -19884.5
